# Update Price Impact (%), Incremental IL (%), and IL/Price Impact (%) columns
# for rows 3-27 to use a forward-looking (next-row) comparison instead of the
# previous backward-looking (prior-row) comparison, matching the refreshed
# simulation output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row, Price Impact (%), Incremental IL (%), IL/Price Impact (%)
$updates = @(
    @(3, -0.0004081820073276177, -0.0004040162853602425, 98.97944497979019),
    @(4, -0.0004081803412270268, -0.0003956832757334716, 96.93834704141118),
    @(5, -0.0004081786751042316, -0.0003873504153428797, 94.8972690070021),
    @(6, -0.000408177009014743, -0.0003790177020901453, 92.85621034977389),
    @(7, -0.0004081753429585611, -0.0003706851338769468, 90.81517055638993),
    @(8, -0.0004081736768690725, -0.0003623527086715761, 88.77414914431287),
    @(9, -0.0004081720108239928, -0.0003540204243090983, 86.73314556635655),
    @(10, -0.0004081703447900153, -0.0003456882788022142, 84.69215934343656),
    @(11, -0.0004081686787782424, -0.0003373562699637844, 82.65118993784178),
    @(12, -0.0004081670127553672, -0.0003290243957621009, 80.61023685892519),
    @(13, -0.0004081653467546964, -0.0003206926541432509, 78.5692995971028),
    @(14, -0.0004081636807873323, -0.0003123610429978108, 76.52837763401146),
    @(15, -0.000408162014808866, -0.0003040295601941523, 74.48747045619204),
    @(16, -0.0004081603488526042, -0.0002956982036450562, 72.4465775463749),
    @(17, -0.0004081586829296491, -0.0002873669714076321, 70.40569842713923),
    @(18, -0.0004081570169844895, -0.0002790358612392296, 68.36483255899368),
    @(19, -0.0004081553510615343, -0.0002707048710970383, 66.32397943405287),
    @(20, -0.0004081536851829881, -0.0002623739989826568, 64.28313855968892),
    @(21, -0.0004081520192711352, -0.0002540432426978434, 62.24230940998544),
    @(22, -0.0004081503534036912, -0.00024571260021089, 60.20149147533921),
    @(23, -0.0004081486875362472, -0.0002373820694345774, 58.16068425149493),
    @(24, -0.0004081470217132122, -0.0002290516483149929, 56.11988722924895),
    @(25, -0.0004081453558568704, -0.0002207213346983039, 54.07909989197748),
    @(26, -0.0004081436900449376, -0.0002123911265750067, 52.03832173703873),
    @(27, 0, 0, 0)
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 6).Value = $u[1]   # F: Price Impact (%)
    $ws.Cells.Item($r, 7).Value = $u[2]   # G: Incremental IL (%)
    $ws.Cells.Item($r, 8).Value = $u[3]   # H: IL/Price Impact (%)
}

